# Split ISIC code 05T06 into ISIC 05 and ISIC 06 on the SoCaOMSbRIC sheet.
#
# The "ISIC 05T06" column (column C) is split into two adjacent columns:
#   C -> "ISIC 05"
#   D -> "ISIC 06" (newly inserted, inheriting column C's original value)
# Every column to the right of the old column C shifts one column to the
# right as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCaOMSbRIC")

# Insert a new column at D, pushing the former "ISIC 07T08" (and everything
# after it) one column to the right. Column C ("ISIC 05T06") keeps its
# position and data, and we relabel C/D as the two new split codes.
$ws.Columns("D").Insert()

$ws.Range("C1").Value = "ISIC 05"
$ws.Range("D1").Value = "ISIC 06"

# The original "ISIC 05T06" data value (0) is preserved in C2; mirror it
# into the newly created D2 cell.
$ws.Range("D2").Value = $ws.Range("C2").Value2
